$wb = $excel.ActiveWorkbook

$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# Sheet "About"
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = "Version: $newVersion"
$wsAbout.Range("A6").Value = "Recommended Citation:  " + '"' + "Global Energy Monitor, Coal mine boundaries and methane sources for Yerunakovskaya-8 Coal Mine, Russia, M0877, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# Sheet "Boundaries and methane sources"
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 11; $row++) {
    $wsData.Cells.Item($row, 19).Value = $newVersion
}
